# Generate Report for Archive
# Update the "Status" cells for the a6db3d67... and bf0cbe50... rows
# from "Ready for handoff" to "In Translation" across the Overview,
# zh-cn, and de-de worksheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns B (zh-cn) and C (de-de) hold the status ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B8").Value = "In Translation"
$wsOverview.Range("C8").Value = "In Translation"
$wsOverview.Range("B9").Value = "In Translation"
$wsOverview.Range("C9").Value = "In Translation"

# --- zh-cn sheet: column C holds Status ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C8").Value = "In Translation"
$wsZh.Range("C9").Value = "In Translation"

# --- de-de sheet: column C holds Status ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C8").Value = "In Translation"
$wsDe.Range("C9").Value = "In Translation"
